# Weekly price-sheet update for "Hortaliza, Agrícola del Norte S.A. de Arica - Betarraga".
#
# A new weekly record (pair of rows: "Primera" and "Segunda") is inserted at the
# top of the historical data block (row 364), pushing the existing data (rows
# 364-470) down by two rows (to 366-472). The two new rows are populated with a
# new date and updated price figures; everything else is carried over unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at row 364 (shifts old rows 364:470 down to 366:472).
$ws.Rows.Item(364).Insert()
$ws.Rows.Item(364).Insert()

# Seed the two new rows with the data that is now one (364) and two (365) rows
# below them (the rows that used to be 364/365 before the insert), so all of
# the non-numeric / unchanged columns (A,B,C,E,F,G,H,I,N,O,Q,R) come along for
# free. The price/date columns are then overwritten with the new week's values.
$ws.Range("A366:R366").Copy($ws.Range("A364:R364"))
$ws.Range("A367:R367").Copy($ws.Range("A365:R365"))

$newDate = Get-Date -Year 2023 -Month 6 -Day 16 -Hour 0 -Minute 0 -Second 0

# Row 364 ("Primera")
$ws.Cells.Item(364, 4).Value  = $newDate   # D: Fecha
$ws.Cells.Item(364, 10).Value = 750         # J
$ws.Cells.Item(364, 12).Value = 700         # L
$ws.Cells.Item(364, 13).Value = 640         # M
$ws.Cells.Item(364, 16).Value = 160         # P

# Row 365 ("Segunda")
$ws.Cells.Item(365, 4).Value  = $newDate   # D: Fecha
$ws.Cells.Item(365, 12).Value = 700         # L
$ws.Cells.Item(365, 13).Value = 650         # M
$ws.Cells.Item(365, 16).Value = 130         # P
